$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the transferred_at date cells from text to real Excel date
# serial values and format them as dates (yyyy-mm-dd), so the
# import/export round-trips dates properly.
$ws.Range("A2").Value = 43831
$ws.Range("A2").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("A3").Value = 43831
$ws.Range("A3").NumberFormat = "yyyy\-mm\-dd"

# Preserve the originally-recorded active cell/selection.
$ws.Range("D25").Select()
